$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated Price values look like plain decimal numbers (e.g. "301.45").
# Excel would normally auto-convert such strings to numeric cells, but the source
# data stores them as text, so we force a Text number format before assigning,
# then restore the default (Normal) style afterwards so formatting is unaffected.
$textForceAddrs = @("D5", "D6", "D7", "D9", "D10", "D11", "D13", "D15", "D16", "D18", "D20", "D21", "D22", "D23", "D25", "D27", "D28", "D29", "D31", "D34", "D36", "D37", "D39", "D40", "D42", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $textForceAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row-by-row value updates
$ws.Range("D2").Value = "42.595.42"
$ws.Range("E2").Value = "  -0.54%  "

$ws.Range("D3").Value = "2.548.38"

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "301.45"
$ws.Range("E5").Value = "  +1.56%  "

$ws.Range("D6").Value = "97.35"
$ws.Range("E6").Value = "  +6.28%  "

$ws.Range("D7").Value = "0.572"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "0.543"
$ws.Range("E9").Value = "  -0.71%  "

$ws.Range("D10").Value = "36.45"
$ws.Range("E10").Value = "  +2.32%  "

$ws.Range("D11").Value = "0.0804"
$ws.Range("E11").Value = "  -0.37%  "

$ws.Range("E12").Value = "  +8.60%  "

$ws.Range("D13").Value = "7.62"
$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("D14").Value = "2.480.06"
$ws.Range("E14").Value = "  -2.35%  "

$ws.Range("D15").Value = "0.872"
$ws.Range("E15").Value = "  +1.12%  "

$ws.Range("D16").Value = "14.58"
$ws.Range("E16").Value = "  +3.32%  "

$ws.Range("D17").Value = "42.621.81"
$ws.Range("E17").Value = "  -0.48%  "

$ws.Range("D18").Value = "13.22"
$ws.Range("E18").Value = "  +6.23%  "

$ws.Range("D19").Value = "0.0₃0980"
$ws.Range("E19").Value = "  +0.53%  "

$ws.Range("D20").Value = "6.55"
$ws.Range("E20").Value = "  -1.26%  "

$ws.Range("D21").Value = "71.41"
$ws.Range("E21").Value = "  -1.26%  "

$ws.Range("D22").Value = "253.61"
$ws.Range("E22").Value = "  -2.55%  "

$ws.Range("D23").Value = "2.93"
$ws.Range("E23").Value = "  +1.84%  "

$ws.Range("E24").Value = "  -2.03%  "

$ws.Range("D25").Value = "27.70"
$ws.Range("E25").Value = "  -6.12%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("D27").Value = "9.98"
$ws.Range("E27").Value = "  -0.14%  "

$ws.Range("D28").Value = "37.84"
$ws.Range("E28").Value = "  +4.95%  "

$ws.Range("D29").Value = "2.09"
$ws.Range("E29").Value = "  -0.93%  "

$ws.Range("E30").Value = "  +0.62%  "

$ws.Range("D31").Value = "155.48"
$ws.Range("E31").Value = "  +3.19%  "

$ws.Range("E32").Value = "  +0.65%  "

$ws.Range("E33").Value = "  +1.36%  "

$ws.Range("D34").Value = "0.0797"
$ws.Range("E34").Value = "  +0.74%  "

$ws.Range("E35").Value = "  -2.93%  "

$ws.Range("D36").Value = "18.29"
$ws.Range("E36").Value = "  +12.91%  "

$ws.Range("D37").Value = "25.69"
$ws.Range("E37").Value = "  +6.25%  "

$ws.Range("E38").Value = "  +0.94%  "

$ws.Range("D39").Value = "0.118"
$ws.Range("E39").Value = "  -0.70%  "

$ws.Range("D40").Value = "2.09"
$ws.Range("E40").Value = "  +32.20%  "

$ws.Range("E41").Value = "  +0.80%  "

$ws.Range("D42").Value = "3.34"
$ws.Range("E42").Value = "  -1.88%  "

$ws.Range("D43").Value = "2.071.42"
$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.13%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0302"
$ws.Range("E45").Value = "  -2.29%  "

$ws.Range("D46").Value = "87.93"
$ws.Range("E46").Value = "  +3.66%  "

$ws.Range("D47").Value = "9.19"
$ws.Range("E47").Value = "  +6.45%  "

$ws.Range("D48").Value = "2.799.14"
$ws.Range("E48").Value = "  +0.28%  "

$ws.Range("D49").Value = "74.40"
$ws.Range("E49").Value = "  +7.70%  "

$ws.Range("D50").Value = "102.76"
$ws.Range("E50").Value = "  -0.97%  "

$ws.Range("D51").Value = "0.187"
$ws.Range("E51").Value = "  +0.96%  "

# Restore default (Normal) style for the forced-text cells so their formatting
# matches the original (unstyled) cells.
foreach ($addr in $textForceAddrs) {
    $ws.Range($addr).Style = "Normal"
}
